# Daily attendance processing - 2025-12-07 12:41:45
# Applies cell-level corrections to the CVS session analysis sheet:
#  - Class Statistics / summary numbers (Missing, Pending, Coverage %, Avg Attendance %)
#  - "Recorded By" email lists re-ordered/updated per latest sync
#  - Two sessions (rows 47 and 85) flip from "Pending" to "Not Recorded" and are
#    re-colored (pink) to match the other "Not Recorded" rows, e.g. row 7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# --- Row 7 ---
$ws.Range("L7").Value = 15

# --- Row 8 ---
$ws.Range("L8").Value = 117

# --- Row 10 ---
$ws.Range("G10").Value = "Madeha.Saeed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg"
$ws.Range("L10").Value = "33.9%"

# --- Row 15 ---
$ws.Range("G15").Value = "nancy.abdelshafy@med.asu.edu.eg, hana.amr@med.asu.edu.eg"
$ws.Range("H15").Value = "108/216"
$ws.Range("S15").Value = "32.7%"

# --- Row 16 ---
$ws.Range("S16").Value = "30.8%"

# --- Row 17 ---
$ws.Range("P17").Value = 3
$ws.Range("Q17").Value = 14

# --- Row 18 ---
$ws.Range("G18").Value = "shorokmohamed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"
$ws.Range("P18").Value = 2
$ws.Range("Q18").Value = 15

# --- Row 19 ---
$ws.Range("G19").Value = "neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"

# --- Row 32 ---
$ws.Range("G32").Value = "Madeha.Saeed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg"

# --- Row 37 ---
$ws.Range("G37").Value = "nancy.abdelshafy@med.asu.edu.eg, hana.amr@med.asu.edu.eg"
$ws.Range("H37").Value = "23/217"

# --- Row 40 ---
$ws.Range("G40").Value = "shorokmohamed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"

# --- Row 41 ---
$ws.Range("G41").Value = "neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"

# --- Row 46 ---
$ws.Range("G46").Value = "nourhan.mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"

# --- Row 47 ---
$ws.Range("I47").Value = "Not Recorded"
# Re-color row 47 (A47:I47) to match the "Not Recorded" style used elsewhere (e.g. row 7)
$ws.Range("A7:I7").Copy() | Out-Null
$ws.Range("A47:I47").PasteSpecial(-4122) | Out-Null

# --- Row 52 ---
$ws.Range("G52").Value = "Shimaa.ashraf@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"

# --- Row 54 ---
$ws.Range("G54").Value = "Madeha.Saeed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, merna.said@med.asu.edu.eg, maimustafa@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg"

# --- Row 62 ---
$ws.Range("G62").Value = "shorokmohamed@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"

# --- Row 68 ---
$ws.Range("G68").Value = "nourhan.mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"

# --- Row 74 ---
$ws.Range("G74").Value = "Shimaa.ashraf@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"

# --- Row 76 ---
$ws.Range("G76").Value = "Madeha.Saeed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, merna.said@med.asu.edu.eg, maimustafa@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg"

# --- Row 84 ---
$ws.Range("G84").Value = "shorokmohamed@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"

# --- Row 85 ---
$ws.Range("I85").Value = "Not Recorded"
# Re-color row 85 (A85:I85) to match the "Not Recorded" style used elsewhere (e.g. row 7)
$ws.Range("A7:I7").Copy() | Out-Null
$ws.Range("A85:I85").PasteSpecial(-4122) | Out-Null

# --- Row 90 ---
$ws.Range("G90").Value = "shaimaa.ahmed@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"

# --- Row 96 ---
$ws.Range("G96").Value = "norhan.mohamed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg"

# --- Row 98 ---
$ws.Range("G98").Value = "Madeha.Saeed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, merna.said@med.asu.edu.eg, maimustafa@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg"

# --- Row 106 ---
$ws.Range("G106").Value = "wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"

# --- Row 107 ---
$ws.Range("G107").Value = "wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"

# --- Row 112 ---
$ws.Range("G112").Value = "shaimaa.ahmed@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"

# --- Row 118 ---
$ws.Range("G118").Value = "norhan.mohamed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg"

# --- Row 120 ---
$ws.Range("G120").Value = "Madeha.Saeed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, merna.said@med.asu.edu.eg, maimustafa@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg"

# --- Row 128 ---
$ws.Range("G128").Value = "wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"

# --- Row 129 ---
$ws.Range("G129").Value = "wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg"

# --- Row 134 ---
$ws.Range("G134").Value = "Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"

# --- Row 142 ---
$ws.Range("G142").Value = "yassmina.fattoh@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, merna.said@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg"

# --- Row 150 ---
$ws.Range("G150").Value = "neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"

# --- Row 156 ---
$ws.Range("G156").Value = "alshimaa.atef@med.asu.edu.egm, Mohammedeltanany@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"

# --- Row 164 ---
$ws.Range("G164").Value = "yassmina.fattoh@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, merna.said@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg"

# --- Row 172 ---
$ws.Range("G172").Value = "neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"

$excel.CutCopyMode = 0

